$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting existing rows 94:134 down to 95:135.
$ws.Rows("94:94").Insert()

# Populate the newly inserted row 94 with this week's new record.
$ws.Range("A94").Value = 4
$ws.Range("B94").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C94").Value = "Los Lagos"
$ws.Range("D94").Value = 44468
$ws.Range("E94").Value = 10
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100104
$ws.Range("H94").Value = "Frutos de pepita"
$ws.Range("I94").Value = 100104005
$ws.Range("J94").Value = "Pera"
$ws.Range("K94").Value = "Packham's Triumph"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 200
$ws.Range("N94").Value = 16000
$ws.Range("O94").Value = 16000
$ws.Range("P94").Value = 16000
$ws.Range("Q94").Value = "$/caja 15 kilos empedrada"
$ws.Range("R94").Value = "Región de O'Higgins"
$ws.Range("S94").Value = 1067
$ws.Range("T94").Value = 15
